$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) study sheet: update the example value for the "workflow" field
# ---------------------------------------------------------------------------
$studyWs = $wb.Worksheets.Item("study")
$studyWs.Unprotect()
$studyWs.Range("C3").Value2 = "e.g. Laser microdissection"
$studyWs.Protect("B881")

# ---------------------------------------------------------------------------
# 2) HiddenDropdowns sheet: expand the "tissue" dropdown list (col F) with
#    30 new terms, keeping the full list alphabetically sorted.
# ---------------------------------------------------------------------------
$tissueValues = @(
    "Abdomen",
    "Anterior Body",
    "Blade",
    "Blood",
    "Bodywall",
    "Bone Marrow Hematopoietic Niches",
    "Bract",
    "Brain",
    "Bud",
    "Cap",
    "Cephalothorax",
    "Cortex Development In Roots",
    "Developing Brain",
    "Developing Embryo",
    "Developing Seed Coat",
    "Dna Extract",
    "Egg",
    "Eggshell",
    "Embryonic Stem Cells",
    "Endocrine Tissue",
    "Endosperm",
    "Eye",
    "Fat Body",
    "Fetal Heart",
    "Fetal Kidney",
    "Fetal Liver",
    "Fetal Lung",
    "Fin",
    "Floral Meristem",
    "Flower",
    "Gastrulation-Stage Embryo",
    "Gill Animal",
    "Gill Fungi",
    "Gonad",
    "Gut",
    "Hair",
    "Head",
    "Heart",
    "Hepatopancreas",
    "Holdfast Fungi",
    "Hypocotyl",
    "Inflorescence",
    "Intestine",
    "Kidney",
    "Leaf",
    "Leaf Primordia",
    "Leg",
    "Limb Bud",
    "Liver",
    "Lung",
    "Mid Body",
    "Modular Colony",
    "Mollusc Foot",
    "Multicellular Organisms In Culture",
    "Muscle",
    "Mycelium",
    "Mycorrhiza",
    "Neonatal Retina",
    "Neural Crest Cells",
    "Nodules",
    "Not Applicable",
    "Not Collected",
    "Not Provided",
    "Olfactory Epithelium",
    "Other Fungal Tissue",
    "Other Plant Tissue",
    "Other Reproductive Animal Tissue",
    "Other Somatic Animal Tissue",
    "Ovary Animal",
    "Oviduct",
    "Ovule Primordium",
    "Pancreas",
    "Petiole",
    "Placental Trophoblast Cells",
    "Pollen Mother Cells",
    "Posterior Body",
    "Postnatal Thymus",
    "Rhizome Meristem ",
    "Root",
    "Root Apical Meristem",
    "Scales",
    "Scat",
    "Seed",
    "Seedling",
    "Shoot",
    "Shoot Apical Meristem",
    "Skin",
    "Somitic Mesoderm",
    "Sperm Seminal Fluid",
    "Spleen",
    "Spore",
    "Spore Bearing Structure",
    "Stem",
    "Stipe",
    "Stomach",
    "Tentacle",
    "Terminal Body",
    "Testis",
    "Thallus Fungi",
    "Thallus Plant",
    "Thorax",
    "Trichome Precursor Cells",
    "Unicellular Organisms In Culture",
    "Vascular Cambium",
    "Whole Organism",
    "Whole Plant"
)

$hiddenWs = $wb.Worksheets.Item("HiddenDropdowns")
for ($i = 0; $i -lt $tissueValues.Length; $i++) {
    $hiddenWs.Cells.Item(5 + $i, 6).Value2 = $tissueValues[$i]
}

# ---------------------------------------------------------------------------
# 3) imaging_protocol sheet: rename section_thickness_method ->
#    section_thickness_measurement_method (header of col G)
# ---------------------------------------------------------------------------
$imagingWs = $wb.Worksheets.Item("imaging_protocol")
$imagingWs.Unprotect()
$imagingWs.Range("G1").Value2 = "section_thickness_measurement_method (optional)"
$imagingWs.Protect("DD04")

# ---------------------------------------------------------------------------
# 4) sample sheet: point the "tissue" column's data validation list at the
#    now-larger HiddenDropdowns range (F5:F110 instead of F5:F80)
# ---------------------------------------------------------------------------
$sampleWs = $wb.Worksheets.Item("sample")
$sampleWs.Unprotect()
$tissueRange = $sampleWs.Range("F5:F1005")
$tissueRange.Validation.Delete()
$tissueRange.Validation.Add(3, 1, 1, "HiddenDropdowns!`$F`$5:`$F`$110")
$tissueRange.Validation.IgnoreBlank = $true
$tissueRange.Validation.InCellDropdown = $true
$tissueRange.Validation.ShowInput = $true
$tissueRange.Validation.ShowError = $true
$tissueRange.Validation.InputMessage = "Choose from the list"
$sampleWs.Protect("F02E")

Write-Host "Edit complete"
